# Update the class schedule sheet: remove some classes, add new ones,
# and re-merge the time-span cells to match the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First, un-merge every merged range that is changing so the cells
#     underneath are free to be edited individually. ---
$ws.Range("F8:F9").UnMerge()
$ws.Range("B12:B13").UnMerge()
$ws.Range("F12:F13").UnMerge()
$ws.Range("C11:C13").UnMerge()
$ws.Range("D8:D9").UnMerge()
$ws.Range("B8:B9").UnMerge()
$ws.Range("E11:E13").UnMerge()
$ws.Range("E6:E7").UnMerge()
$ws.Range("C6:C7").UnMerge()
$ws.Range("D12:D13").UnMerge()

# --- Row 4: add CPSC_V 221-L1K ---
$ws.Range("D4").Value = "CPSC_V 221-L1K - Basic Algorithms and Data Structures`n9:00 a.m. - 11:00 a.m."

# --- Row 5: add ENGL_V 111 sections ---
$ws.Range("C5").Value = "ENGL_V 111-002 - Approaches to Language and Communication`n9:30 a.m. - 11:00 a.m."
$ws.Range("E5").Value = "ENGL_V 111-L10 - Approaches to Language and Communication`n9:30 a.m. - 11:00 a.m."

# --- Row 6: remove CPSC_V 213-L2F ---
$ws.Range("C6").Value = ""
$ws.Range("E6").Value = ""

# --- Row 8: remove CPSC_V 213-205, add MATH_V 200-102 ---
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = "MATH_V 200-102 - Calculus III`n11:00 a.m. - 12:30 p.m."
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = "MATH_V 200-102 - Calculus III`n11:00 a.m. - 12:30 p.m."
$ws.Range("F8").Value = ""

# --- Row 11: remove MATH_V 221-202 ---
$ws.Range("C11").Value = ""
$ws.Range("E11").Value = ""

# --- Row 12: remove ENGL_V 110 sections, add CPSC_V 330-T1F ---
$ws.Range("B12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = "CPSC_V 330-T1F - Applied Machine Learning`n1:00 p.m. - 2:00 p.m."
$ws.Range("F12").Value = ""

# --- Row 14: add CPSC_V 221-101 ---
$ws.Range("B14").Value = "CPSC_V 221-101 - Basic Algorithms and Data Structures`n2:00 p.m. - 3:00 p.m."
$ws.Range("D14").Value = "CPSC_V 221-101 - Basic Algorithms and Data Structures`n2:00 p.m. - 3:00 p.m."
$ws.Range("F14").Value = "CPSC_V 221-101 - Basic Algorithms and Data Structures`n2:00 p.m. - 3:00 p.m."

# --- Row 17: add CPSC_V 330-101 ---
$ws.Range("C17").Value = "CPSC_V 330-101 - Applied Machine Learning`n3:30 p.m. - 5:00 p.m."
$ws.Range("E17").Value = "CPSC_V 330-101 - Applied Machine Learning`n3:30 p.m. - 5:00 p.m."

# --- Re-create the merged ranges to match the new layout ---
$ws.Range("D4:D7").Merge()
$ws.Range("E17:E19").Merge()
$ws.Range("E8:E10").Merge()
$ws.Range("E12:E13").Merge()
$ws.Range("C17:C19").Merge()
$ws.Range("C5:C7").Merge()
$ws.Range("E5:E7").Merge()
$ws.Range("C8:C10").Merge()
$ws.Range("D14:D15").Merge()
$ws.Range("B14:B15").Merge()
$ws.Range("F14:F15").Merge()
